$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add new date header in column C, matching B1 style (bold/centered/bordered)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rows 2-42: reordered fund rows (funds shift up by 2, total/avg move to bottom)
# Column A = label, Column B = unchanged prior value, Column C = new value (13-01-2023)
$ws.Cells.Item(2, 1).Value = "1810 Renta variable"
$ws.Cells.Item(2, 2).Value = 228540.77
$ws.Cells.Item(2, 3).Value = 228317.61
$ws.Cells.Item(3, 1).Value = "1822 Raices Valores Negociables"
$ws.Cells.Item(3, 2).Value = 151779.56
$ws.Cells.Item(3, 3).Value = 150021.85
$ws.Cells.Item(4, 1).Value = "Adcap IOL Acciones Argentina"
$ws.Cells.Item(4, 2).Value = 41212.78
$ws.Cells.Item(4, 3).Value = 41305.64
$ws.Cells.Item(5, 1).Value = "Alpha Mega"
$ws.Cells.Item(5, 2).Value = 180966.68
$ws.Cells.Item(5, 3).Value = 180612.28
$ws.Cells.Item(6, 1).Value = "Alpha planeam equil"
$ws.Cells.Item(6, 2).Value = 8938.06
$ws.Cells.Item(6, 3).Value = 8963.26
$ws.Cells.Item(7, 1).Value = "Alpha renta balan global"
$ws.Cells.Item(7, 2).Value = 667763.02
$ws.Cells.Item(7, 3).Value = 668178.99
$ws.Cells.Item(8, 1).Value = "Argenfunds"
$ws.Cells.Item(8, 2).Value = 9787.89
$ws.Cells.Item(8, 3).Value = 9774.53
$ws.Cells.Item(9, 1).Value = "Arpenta acciones"
$ws.Cells.Item(9, 2).Value = 15167.95
$ws.Cells.Item(9, 3).Value = 15159.23
$ws.Cells.Item(10, 1).Value = "Arpenta ex Mercosur"
$ws.Cells.Item(10, 2).Value = 44972.96
$ws.Cells.Item(10, 3).Value = 45025.72
$ws.Cells.Item(11, 1).Value = "Balanz"
$ws.Cells.Item(11, 2).Value = 230141.76
$ws.Cells.Item(11, 3).Value = 231091.38
$ws.Cells.Item(12, 1).Value = "Compass Crecimiento"
$ws.Cells.Item(12, 2).Value = 1588335.52
$ws.Cells.Item(12, 3).Value = 1579787.3
$ws.Cells.Item(13, 1).Value = "Consultatio Acciones Argentina"
$ws.Cells.Item(13, 2).Value = 241153.55
$ws.Cells.Item(13, 3).Value = 240434.53
$ws.Cells.Item(14, 1).Value = "Consultatio Renta Variable"
$ws.Cells.Item(14, 2).Value = 234157.45
$ws.Cells.Item(14, 3).Value = 234177.03
$ws.Cells.Item(15, 1).Value = "Delta Select"
$ws.Cells.Item(15, 2).Value = 345053.22
$ws.Cells.Item(15, 3).Value = 318887.24
$ws.Cells.Item(16, 1).Value = "Delta gestion V"
$ws.Cells.Item(16, 2).Value = 92115.71
$ws.Cells.Item(16, 3).Value = 92103.73
$ws.Cells.Item(17, 1).Value = "FBA Acciones Argentinas"
$ws.Cells.Item(17, 2).Value = 100357
$ws.Cells.Item(17, 3).Value = 100532.67
$ws.Cells.Item(18, 1).Value = "FBA Calificado"
$ws.Cells.Item(18, 2).Value = 98716.33
$ws.Cells.Item(18, 3).Value = 99881.43
$ws.Cells.Item(19, 1).Value = "Fima Acciones"
$ws.Cells.Item(19, 2).Value = 667794.67
$ws.Cells.Item(19, 3).Value = 666876.98
$ws.Cells.Item(20, 1).Value = "Fima PB Acciones"
$ws.Cells.Item(20, 2).Value = 170767.56
$ws.Cells.Item(20, 3).Value = 171096.47
$ws.Cells.Item(21, 1).Value = "Galileo Acciones"
$ws.Cells.Item(21, 2).Value = 3923342
$ws.Cells.Item(21, 3).Value = 3924173.45
$ws.Cells.Item(22, 1).Value = "Goal Acciones Argentinas"
$ws.Cells.Item(22, 2).Value = 44513.71
$ws.Cells.Item(22, 3).Value = 44699.71
$ws.Cells.Item(23, 1).Value = "Goal acciones plus"
$ws.Cells.Item(23, 2).Value = 2015.97
$ws.Cells.Item(23, 3).Value = 2001.29
$ws.Cells.Item(24, 1).Value = "HF Acciones Argentinas"
$ws.Cells.Item(24, 2).Value = 131464.87
$ws.Cells.Item(24, 3).Value = 131660.55
$ws.Cells.Item(25, 1).Value = "HF Acciones Lideres"
$ws.Cells.Item(25, 2).Value = 179125.87
$ws.Cells.Item(25, 3).Value = 177945.73
$ws.Cells.Item(26, 1).Value = "IAM Renta Variable"
$ws.Cells.Item(26, 2).Value = 31611.5
$ws.Cells.Item(26, 3).Value = 32867.42
$ws.Cells.Item(27, 1).Value = "IEB Value"
$ws.Cells.Item(27, 2).Value = 3172.93
$ws.Cells.Item(27, 3).Value = 3195.36
$ws.Cells.Item(28, 1).Value = "Lombardi"
$ws.Cells.Item(28, 2).Value = 17221.56
$ws.Cells.Item(28, 3).Value = 17199.12
$ws.Cells.Item(29, 1).Value = "MAF"
$ws.Cells.Item(29, 2).Value = 17914.92
$ws.Cells.Item(29, 3).Value = 18105.08
$ws.Cells.Item(30, 1).Value = "Megainver"
$ws.Cells.Item(30, 2).Value = 27070.97
$ws.Cells.Item(30, 3).Value = 27054.06
$ws.Cells.Item(31, 1).Value = "Pellegrini Acciones"
$ws.Cells.Item(31, 2).Value = 98230.83
$ws.Cells.Item(31, 3).Value = 98039.91
$ws.Cells.Item(32, 1).Value = "Pionero Acciones"
$ws.Cells.Item(32, 2).Value = 59223
$ws.Cells.Item(32, 3).Value = 59055.94
$ws.Cells.Item(33, 1).Value = "Premier Renta Variable"
$ws.Cells.Item(33, 2).Value = 70754.56
$ws.Cells.Item(33, 3).Value = 70741.14
$ws.Cells.Item(34, 1).Value = "Quinquela Acciones"
$ws.Cells.Item(34, 2).Value = 86322.5
$ws.Cells.Item(34, 3).Value = 86828.71
$ws.Cells.Item(35, 1).Value = "Rofex 20 Renta Variable"
$ws.Cells.Item(35, 2).Value = 62224.84
$ws.Cells.Item(35, 3).Value = 62146.73
$ws.Cells.Item(36, 1).Value = "SBS Acciones Argentina"
$ws.Cells.Item(36, 2).Value = 141434.49
$ws.Cells.Item(36, 3).Value = 238036.54
$ws.Cells.Item(37, 1).Value = "Schroeder RV"
$ws.Cells.Item(37, 2).Value = 376145.49
$ws.Cells.Item(37, 3).Value = 376234.52
$ws.Cells.Item(38, 1).Value = "Supefondo RV"
$ws.Cells.Item(38, 2).Value = 1376831.53
$ws.Cells.Item(38, 3).Value = 1379923.04
$ws.Cells.Item(39, 1).Value = "Superfondo "
$ws.Cells.Item(39, 2).Value = 2091740.24
$ws.Cells.Item(39, 3).Value = 2092740.93
$ws.Cells.Item(40, 1).Value = "Toronto Trust Multimercado"
$ws.Cells.Item(40, 2).Value = 40058.31
$ws.Cells.Item(40, 3).Value = 40239.96
$ws.Cells.Item(41, 1).Value = "avg"
$ws.Cells.Item(41, 2).Value = 356362.63
$ws.Cells.Item(41, 3).Value = 358079.92
$ws.Cells.Item(42, 1).Value = "total"
$ws.Cells.Item(42, 2).Value = 13898142.53
$ws.Cells.Item(42, 3).Value = 13965117.06
